# Update "Resumo de Inscricoes" counts (Inscritos / Pagos / Inscricoes homologadas)
# for several Processo/Cargo rows on the "Inscricoes" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 11
$ws.Range("E11").Value = 287
$ws.Range("F11").Value = 160
$ws.Range("H11").Value = 160

# Row 12
$ws.Range("E12").Value = 414

# Row 13
$ws.Range("E13").Value = 109

# Row 21
$ws.Range("E21").Value = 126

# Row 25
$ws.Range("E25").Value = 223
$ws.Range("F25").Value = 103
$ws.Range("H25").Value = 103

# Row 27
$ws.Range("E27").Value = 284
$ws.Range("F27").Value = 135
$ws.Range("H27").Value = 135

# Row 28
$ws.Range("E28").Value = 168
$ws.Range("F28").Value = 58
$ws.Range("H28").Value = 58

# Row 29
$ws.Range("E29").Value = 148

# Row 31
$ws.Range("E31").Value = 67
$ws.Range("F31").Value = 31
$ws.Range("H31").Value = 31

# Row 33
$ws.Range("E33").Value = 250

# Row 34
$ws.Range("E34").Value = 182

# Row 35
$ws.Range("E35").Value = 122

# Row 36
$ws.Range("E36").Value = 57
$ws.Range("F36").Value = 35
$ws.Range("H36").Value = 35

# Row 37
$ws.Range("E37").Value = 134

# Row 38
$ws.Range("E38").Value = 83

# Row 40
$ws.Range("E40").Value = 230
$ws.Range("F40").Value = 106
$ws.Range("H40").Value = 106

# Row 41
$ws.Range("E41").Value = 337

# Row 42
$ws.Range("E42").Value = 309
$ws.Range("F42").Value = 165
$ws.Range("H42").Value = 165

# Row 45
$ws.Range("E45").Value = 122
$ws.Range("F45").Value = 60
$ws.Range("H45").Value = 60

# Row 46
$ws.Range("E46").Value = 270

# Row 47
$ws.Range("E47").Value = 381

# Row 48
$ws.Range("F48").Value = 73
$ws.Range("H48").Value = 73

# Row 49
$ws.Range("E49").Value = 256

# Row 50
$ws.Range("E50").Value = 220
$ws.Range("F50").Value = 96
$ws.Range("H50").Value = 96

# Row 51
$ws.Range("E51").Value = 208
$ws.Range("F51").Value = 85
$ws.Range("H51").Value = 85
